$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "71÷8=8, 7" "42÷9=4, 6"
Replace-Text "64÷7=9, 1" "26÷5=5, 1"
Replace-Text "83÷5=16, 3" "70÷3=23, 1"
Replace-Text "88÷3=29, 1" "14÷5=2, 4"
Replace-Text "52÷5=10, 2" "78÷9=8, 6"
Replace-Text "86÷6=14, 2" "74÷9=8, 2"
Replace-Text "35÷5=7, 0" "46÷4=11, 2"
Replace-Text "66÷3=22, 0" "60÷6=10, 0"
Replace-Text "27÷2=13, 1" "96÷4=24, 0"
Replace-Text "12÷2=6, 0" "21÷6=3, 3"
Replace-Text "77÷6=12, 5" "70÷9=7, 7"
Replace-Text "88÷5=17, 3" "37÷3=12, 1"
Replace-Text "48÷5=9, 3" "15÷7=2, 1"
Replace-Text "50÷2=25, 0" "68÷5=13, 3"
Replace-Text "55÷3=18, 1" "76÷5=15, 1"
Replace-Text "75÷3=25, 0" "45÷8=5, 5"
Replace-Text "15÷4=3, 3" "67÷5=13, 2"
Replace-Text "99÷4=24, 3" "75÷8=9, 3"
Replace-Text "76÷9=8, 4" "53÷7=7, 4"
Replace-Text "35÷3=11, 2" "69÷9=7, 6"
Replace-Text "26÷9=2, 8" "72÷8=9, 0"
Replace-Text "46÷6=7, 4" "12÷3=4, 0"
Replace-Text "60÷9=6, 6" "87÷2=43, 1"
Replace-Text "20÷6=3, 2" "31÷8=3, 7"
Replace-Text "87÷9=9, 6" "69÷9=7, 6"

Write-Output "Done"
